$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.863.16"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.327.93"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'582.61"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "'176.65"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "3.324.02"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'47.26"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "'696.56"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "3.872.01"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "67.936.66"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.118"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.337.72"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'11.05"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "'0.896"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "'5.39"
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").Value = "'17.05"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "'99.71"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "'3.91"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'9.60"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("D29").Value = "'33.14"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").Value = "'7.10"
$ws.Range("E31").Value = "  +6.10%  "
$ws.Range("D32").Value = "'568.59"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "3.684.43"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("D38").Value = "'3.38"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").Value = "'34.61"
$ws.Range("E39").Value = "  +8.72%  "
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").Value = "'3.17"
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("D43").Value = "'3.34"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0675"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.337"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").Value = "'0.0406"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").Value = "'2.69"
$ws.Range("E47").Value = "  +4.47%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'130.40"
$ws.Range("E51").Value = "  +0.10%  "
